$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source added a new price record for Femacal de La Calera - Haba that
# belongs right after the existing row 137, so insert a fresh row at 138
# (this shifts the former rows 138:201 down to 139:202, matching the diff).
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new record's data
$ws.Cells.Item(138, 1).Value = 3
$ws.Cells.Item(138, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(138, 3).Value = "Coquimbo"
$ws.Cells.Item(138, 4).Value = 44825
$ws.Cells.Item(138, 5).Value = 5
$ws.Cells.Item(138, 6).Value = 100112026
$ws.Cells.Item(138, 7).Value = "Haba"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 85
$ws.Cells.Item(138, 11).Value = 11500
$ws.Cells.Item(138, 12).Value = 12000
$ws.Cells.Item(138, 13).Value = 11735
$ws.Cells.Item(138, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(138, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(138, 16).Value = 469
$ws.Cells.Item(138, 17).Value = 25
$ws.Cells.Item(138, 18).Value = "Hortaliza"
